$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Mon Feb 24 22:49:50 EST 2025",
    "Mon Feb 24 22:50:03 EST 2025",
    "Mon Feb 24 22:50:15 EST 2025",
    "Mon Feb 24 22:50:28 EST 2025",
    "Mon Feb 24 22:50:42 EST 2025",
    "Mon Feb 24 22:50:55 EST 2025",
    "Mon Feb 24 22:51:07 EST 2025",
    "Mon Feb 24 22:51:20 EST 2025",
    "Mon Feb 24 22:51:33 EST 2025",
    "Mon Feb 24 22:51:46 EST 2025",
    "Mon Feb 24 22:51:59 EST 2025",
    "Mon Feb 24 22:52:12 EST 2025",
    "Mon Feb 24 22:52:25 EST 2025",
    "Mon Feb 24 22:52:38 EST 2025",
    "Mon Feb 24 22:52:51 EST 2025",
    "Mon Feb 24 22:53:05 EST 2025",
    "Mon Feb 24 22:53:18 EST 2025",
    "Mon Feb 24 22:53:36 EST 2025",
    "Mon Feb 24 22:53:55 EST 2025",
    "Mon Feb 24 22:54:09 EST 2025",
    "Mon Feb 24 22:54:22 EST 2025",
    "Mon Feb 24 22:54:35 EST 2025",
    "Mon Feb 24 22:54:48 EST 2025",
    "Mon Feb 24 22:55:01 EST 2025",
    "Mon Feb 24 22:55:15 EST 2025",
    "Mon Feb 24 22:55:28 EST 2025",
    "Mon Feb 24 22:55:41 EST 2025",
    "Mon Feb 24 22:55:55 EST 2025",
    "Mon Feb 24 22:56:09 EST 2025",
    "Mon Feb 24 22:56:22 EST 2025",
    "Mon Feb 24 22:56:35 EST 2025",
    "Mon Feb 24 22:56:48 EST 2025",
    "Mon Feb 24 22:57:03 EST 2025",
    "Mon Feb 24 22:57:15 EST 2025",
    "Mon Feb 24 22:57:28 EST 2025",
    "Mon Feb 24 22:57:42 EST 2025",
    "Mon Feb 24 22:57:55 EST 2025",
    "Mon Feb 24 22:58:10 EST 2025",
    "Mon Feb 24 22:58:25 EST 2025",
    "Mon Feb 24 22:58:53 EST 2025",
    "Mon Feb 24 22:59:06 EST 2025",
    "Mon Feb 24 22:59:20 EST 2025",
    "Mon Feb 24 22:59:34 EST 2025",
    "Mon Feb 24 22:59:48 EST 2025",
    "Mon Feb 24 23:00:01 EST 2025",
    "Mon Feb 24 23:00:15 EST 2025",
    "Mon Feb 24 23:00:29 EST 2025",
    "Mon Feb 24 23:00:45 EST 2025",
    "Mon Feb 24 23:01:03 EST 2025",
    "Mon Feb 24 23:01:18 EST 2025",
    "Mon Feb 24 23:01:31 EST 2025",
    "Mon Feb 24 23:01:45 EST 2025",
    "Mon Feb 24 23:02:00 EST 2025",
    "Mon Feb 24 23:02:13 EST 2025",
    "Mon Feb 24 23:02:26 EST 2025",
    "Mon Feb 24 23:02:40 EST 2025",
    "Mon Feb 24 23:02:53 EST 2025",
    "Mon Feb 24 23:03:06 EST 2025"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
